# Generate Report for Handback
# Update the handoff/handback timestamps and the "Latest HO Xliff Generate Date"
# for the row corresponding to the 3c2fe154-... file across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G, row 2
$overview.Range("G2").Value = "2016-09-03 14:52:35"

# zh-cn sheet: "Correspond Handoff Datetime" column H, row 2
$zhcn.Range("H2").Value = "2016-09-03 14:52:31"
# zh-cn sheet: "Correspond Handback DateTime" column K, row 2
$zhcn.Range("K2").Value = "2016-09-03 14:52:47"

# de-de sheet: "Correspond Handoff Datetime" column H, row 2
$dede.Range("H2").Value = "2016-09-03 14:52:35"
# de-de sheet: "Correspond Handback DateTime" column K, row 2
$dede.Range("K2").Value = "2016-09-03 14:52:54"
